$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph. Built via InsertXML so we reproduce the exact
#    run layout used throughout this document: a leading empty run,
#    a bold run for the label, then a plain run for the rest of the text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our review of Cat Wilde and the Eclipse of the Sun God, a Play''n Go slot game with stunning graphics and medium volatility. Play free now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Near the end of the document, the duplicated bold title paragraph
#    ("Play Cat Wilde and the Eclipse of the Sun God for Free") is
#    removed entirely. Locate it by scanning paragraphs for the title
#    text that is NOT the Heading 1 at the very top of the document
#    (Paragraph.Range.Text carries a trailing paragraph-mark character,
#    hence the TrimEnd).
# ---------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Play Cat Wilde and the Eclipse of the Sun God for Free" -and $para.Style.NameLocal -ne "Heading 1") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3) The last paragraph's italic text is replaced with the new prompt
#    text. We rewrite the text directly on the run's range (not via
#    Find/Replace) so that straight quotes survive untouched by any
#    smart-quote autocorrect behaviour.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = "Please create a feature image for `"Cat Wilde and the Eclipse of the Sun God`" fitting the following criteria: - Cartoon-style image - Happy Maya warrior with glasses"
